$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.538.21"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").Value = "1.724.76"
$ws.Range("E3").Value = "  +4.24%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'225.88"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").Value = "'0.5380"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "'0.06610"
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "'21.80"
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'4.619"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.722.28"
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("D14").Value = "1.962.05"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "'0.5881"
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "0.0₅8316"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "'68.06"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "27.559.18"
$ws.Range("E18").Value = "  +5.42%  "
$ws.Range("D19").Value = "'222.06"
$ws.Range("E19").Value = "  +15.45%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'4.748"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "'6.102"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'148.27"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("E26").Value = "  +12.16%  "
$ws.Range("D27").Value = "'0.1233"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("D28").Value = "'7.402"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "'16.68"
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "'0.05536"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "'3.546"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "'3.463"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").Value = "'1.664"
$ws.Range("E34").Value = "  +6.28%  "
$ws.Range("D35").Value = "'0.9621"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.822"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.446"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'0.5958"
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("D39").Value = "'0.01646"
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("D40").Value = "'5.928"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "1.060.25"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "'0.8538"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'101.61"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "1.868.04"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("E46").Value = "  +13.70%  "
$ws.Range("D47").Value = "'59.11"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "'8.225"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "'0.4439"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'0.05279"
$ws.Range("E51").Value = "  +1.68%  "
